$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 22 ("cos") - update raw input values
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 4918
$ws.Range("G22").Value = 4441
$ws.Range("M22").Value = 1803
$ws.Range("N22").Value = 1668
$ws.Range("W22").Value = 2265
$ws.Range("X22").Value = 1679

# Row 23 ("tan") - update raw input values
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 6137
$ws.Range("G23").Value = 5354
$ws.Range("M23").Value = 2358
$ws.Range("N23").Value = 1673
$ws.Range("W23").Value = 2886
$ws.Range("X23").Value = 2224

$excel.CalculateFullRebuild()
